$wb = $excel.ActiveWorkbook

# Rename the "Lockers" sheet to "Basement" so the tab naming is consistent
# with the other floor-level sheets (Second Floor, Third Floor, Fourth Floor)
# now that locker allocations are being assigned and published per-floor.
$ws = $wb.Worksheets.Item("Lockers")
$ws.Name = "Basement"
